# Generate Report for Handoff
# Adds two newly-handed-off files (3885a69e-... and c0a4997a-...) to the
# localization status report, on all three sheets (Overview, zh-cn, de-de).
# The previously-last row (.localization-config) is pushed down to make
# room for them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Duplicate the last data row (row 4, ".localization-config") twice so we
# end up with two fresh rows (4 and 5) above it, and the config row lands
# on row 6 - preserving every cell's style along the way.
$wsOverview.Rows.Item(4).Copy()
$wsOverview.Rows.Item(4).Insert(-4121)
$wsOverview.Rows.Item(4).Copy()
$wsOverview.Rows.Item(4).Insert(-4121)

$wsOverview.Range("A4").Value = "3885a69e-28b9-473e-ba09-6dc80d75e0f7.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

$wsOverview.Range("A5").Value = "c0a4997a-4c20-4cd6-b69c-a60f1030951d.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

# Rebuild hyperlinks (row insert does not re-anchor the existing ones).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/9d56d942-f46d-4bea-9fa0-51b23aa8edec.md", "", "", "9d56d942-f46d-4bea-9fa0-51b23aa8edec.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/b4e3280a-3a2d-4184-99f5-0667a18af705.md", "", "", "b4e3280a-3a2d-4184-99f5-0667a18af705.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/3885a69e-28b9-473e-ba09-6dc80d75e0f7.md", "", "", "3885a69e-28b9-473e-ba09-6dc80d75e0f7.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/c0a4997a-4c20-4cd6-b69c-a60f1030951d.md", "", "", "c0a4997a-4c20-4cd6-b69c-a60f1030951d.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Duplicate row 3 (a fully-populated "In Translation" row, including the
# Latest Handoff File column) twice to create rows 4 and 5; this pushes
# the ".localization-config" row from 4 down to 6.
$wsZh.Rows.Item(3).Copy()
$wsZh.Rows.Item(4).Insert(-4121)
$wsZh.Rows.Item(3).Copy()
$wsZh.Rows.Item(4).Insert(-4121)

$wsZh.Range("A4").Value = "3885a69e-28b9-473e-ba09-6dc80d75e0f7.md"
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = "3885a69e-28b9-473e-ba09-6dc80d75e0f7.a0d87f6f43cf3ba03cd09648509ea7fed301ecaa.zh-cn.xlf"
$wsZh.Range("D4").Value = "2016-03-10 05:32:28"
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"

$wsZh.Range("A5").Value = "c0a4997a-4c20-4cd6-b69c-a60f1030951d.md"
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Range("C5").Value = "c0a4997a-4c20-4cd6-b69c-a60f1030951d.4537473abad44ea7c812f6506c403e975a6123bd.zh-cn.xlf"
$wsZh.Range("D5").Value = "2016-03-10 05:32:28"
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/9d56d942-f46d-4bea-9fa0-51b23aa8edec.md", "", "", "9d56d942-f46d-4bea-9fa0-51b23aa8edec.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85896189ed16554cfccd2938a286ed89f2356413/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/9d56d942-f46d-4bea-9fa0-51b23aa8edec.59f3184b2239203cf9bbcc4db5af692ddd62c11d.zh-cn.xlf", "", "", "9d56d942-f46d-4bea-9fa0-51b23aa8edec.59f3184b2239203cf9bbcc4db5af692ddd62c11d.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/b4e3280a-3a2d-4184-99f5-0667a18af705.md", "", "", "b4e3280a-3a2d-4184-99f5-0667a18af705.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85896189ed16554cfccd2938a286ed89f2356413/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b4e3280a-3a2d-4184-99f5-0667a18af705.f74b0eb73705186f610f696898b1a9bbec28bc18.zh-cn.xlf", "", "", "b4e3280a-3a2d-4184-99f5-0667a18af705.f74b0eb73705186f610f696898b1a9bbec28bc18.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/3885a69e-28b9-473e-ba09-6dc80d75e0f7.md", "", "", "3885a69e-28b9-473e-ba09-6dc80d75e0f7.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85896189ed16554cfccd2938a286ed89f2356413/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/3885a69e-28b9-473e-ba09-6dc80d75e0f7.a0d87f6f43cf3ba03cd09648509ea7fed301ecaa.zh-cn.xlf", "", "", "3885a69e-28b9-473e-ba09-6dc80d75e0f7.a0d87f6f43cf3ba03cd09648509ea7fed301ecaa.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/c0a4997a-4c20-4cd6-b69c-a60f1030951d.md", "", "", "c0a4997a-4c20-4cd6-b69c-a60f1030951d.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85896189ed16554cfccd2938a286ed89f2356413/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c0a4997a-4c20-4cd6-b69c-a60f1030951d.4537473abad44ea7c812f6506c403e975a6123bd.zh-cn.xlf", "", "", "c0a4997a-4c20-4cd6-b69c-a60f1030951d.4537473abad44ea7c812f6506c403e975a6123bd.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(3).Copy()
$wsDe.Rows.Item(4).Insert(-4121)
$wsDe.Rows.Item(3).Copy()
$wsDe.Rows.Item(4).Insert(-4121)

$wsDe.Range("A4").Value = "3885a69e-28b9-473e-ba09-6dc80d75e0f7.md"
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = "3885a69e-28b9-473e-ba09-6dc80d75e0f7.a0d87f6f43cf3ba03cd09648509ea7fed301ecaa.de-de.xlf"
$wsDe.Range("D4").Value = "2016-03-10 05:32:36"
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"

$wsDe.Range("A5").Value = "c0a4997a-4c20-4cd6-b69c-a60f1030951d.md"
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Range("C5").Value = "c0a4997a-4c20-4cd6-b69c-a60f1030951d.4537473abad44ea7c812f6506c403e975a6123bd.de-de.xlf"
$wsDe.Range("D5").Value = "2016-03-10 05:32:36"
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/9d56d942-f46d-4bea-9fa0-51b23aa8edec.md", "", "", "9d56d942-f46d-4bea-9fa0-51b23aa8edec.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b70a64d2d76177d57a77362e384cb44fe68726a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/9d56d942-f46d-4bea-9fa0-51b23aa8edec.59f3184b2239203cf9bbcc4db5af692ddd62c11d.de-de.xlf", "", "", "9d56d942-f46d-4bea-9fa0-51b23aa8edec.59f3184b2239203cf9bbcc4db5af692ddd62c11d.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/b4e3280a-3a2d-4184-99f5-0667a18af705.md", "", "", "b4e3280a-3a2d-4184-99f5-0667a18af705.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b70a64d2d76177d57a77362e384cb44fe68726a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b4e3280a-3a2d-4184-99f5-0667a18af705.f74b0eb73705186f610f696898b1a9bbec28bc18.de-de.xlf", "", "", "b4e3280a-3a2d-4184-99f5-0667a18af705.f74b0eb73705186f610f696898b1a9bbec28bc18.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/3885a69e-28b9-473e-ba09-6dc80d75e0f7.md", "", "", "3885a69e-28b9-473e-ba09-6dc80d75e0f7.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b70a64d2d76177d57a77362e384cb44fe68726a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/3885a69e-28b9-473e-ba09-6dc80d75e0f7.a0d87f6f43cf3ba03cd09648509ea7fed301ecaa.de-de.xlf", "", "", "3885a69e-28b9-473e-ba09-6dc80d75e0f7.a0d87f6f43cf3ba03cd09648509ea7fed301ecaa.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/e2e/c0a4997a-4c20-4cd6-b69c-a60f1030951d.md", "", "", "c0a4997a-4c20-4cd6-b69c-a60f1030951d.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b70a64d2d76177d57a77362e384cb44fe68726a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c0a4997a-4c20-4cd6-b69c-a60f1030951d.4537473abad44ea7c812f6506c403e975a6123bd.de-de.xlf", "", "", "c0a4997a-4c20-4cd6-b69c-a60f1030951d.4537473abad44ea7c812f6506c403e975a6123bd.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/684ba978ea5fdf88ea48163f5cb3bea07295a811/.localization-config", "", "", ".localization-config")

Write-Output "Applied handoff rows to Overview, zh-cn, de-de sheets."
